# Fruta / hortaliza, semanal
# Insert two new weekly price records (Vega Monumental Concepción - Frutilla)
# right before the existing row 83 (chronologically they are the most recent
# records, dated 2021-10-13), pushing all subsequent rows down by two
# positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 83. Excel shifts rows 83-159
# down to rows 85-161 and the sheet's dimension/UsedRange grows accordingly.
$ws.Rows.Item(83).Insert()
$ws.Rows.Item(83).Insert()

# --- New row 83: "Especial" quality record ---
$ws.Cells.Item(83,1).Value()  = 11
$ws.Cells.Item(83,2).Value()  = "Vega Monumental Concepción"
$ws.Cells.Item(83,3).Value()  = "Bíobío"
$ws.Cells.Item(83,4).Value()  = "2021-10-13"
$ws.Cells.Item(83,5).Value()  = 8
$ws.Cells.Item(83,6).Value()  = "Fruta"
$ws.Cells.Item(83,7).Value()  = 100101
$ws.Cells.Item(83,8).Value()  = "Berries"
$ws.Cells.Item(83,9).Value()  = 100112025
$ws.Cells.Item(83,10).Value() = "Frutilla"
$ws.Cells.Item(83,11).Value() = "Sin especificar"
$ws.Cells.Item(83,12).Value() = "Especial"
$ws.Cells.Item(83,13).Value() = 100
$ws.Cells.Item(83,14).Value() = 12000
$ws.Cells.Item(83,15).Value() = 12000
$ws.Cells.Item(83,16).Value() = 12000
$ws.Cells.Item(83,17).Value() = "`$/bandeja 7 kilos"
$ws.Cells.Item(83,18).Value() = "Provincia de Cardenal Caro"
$ws.Cells.Item(83,19).Value() = 1714
$ws.Cells.Item(83,20).Value() = 7

# --- New row 84: "Primera" quality record ---
$ws.Cells.Item(84,1).Value()  = 11
$ws.Cells.Item(84,2).Value()  = "Vega Monumental Concepción"
$ws.Cells.Item(84,3).Value()  = "Bíobío"
$ws.Cells.Item(84,4).Value()  = "2021-10-13"
$ws.Cells.Item(84,5).Value()  = 8
$ws.Cells.Item(84,6).Value()  = "Fruta"
$ws.Cells.Item(84,7).Value()  = 100101
$ws.Cells.Item(84,8).Value()  = "Berries"
$ws.Cells.Item(84,9).Value()  = 100112025
$ws.Cells.Item(84,10).Value() = "Frutilla"
$ws.Cells.Item(84,11).Value() = "Sin especificar"
$ws.Cells.Item(84,12).Value() = "Primera"
$ws.Cells.Item(84,13).Value() = 100
$ws.Cells.Item(84,14).Value() = 10000
$ws.Cells.Item(84,15).Value() = 10000
$ws.Cells.Item(84,16).Value() = 10000
$ws.Cells.Item(84,17).Value() = "`$/bandeja 7 kilos"
$ws.Cells.Item(84,18).Value() = "Provincia de Cardenal Caro"
$ws.Cells.Item(84,19).Value() = 1429
$ws.Cells.Item(84,20).Value() = 7
